$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Games")
$r = $ws.Range("B94")
foreach ($idx in 5,6,7,8,9,10,11,12) {
  try {
    $b = $r.Borders.Item($idx)
    $b.LineStyle = -4142
  } catch {}
}
Write-Output "done"
